# Insert a new weekly price record for Albahaca (Vega Central Mapocho de
# Santiago) as row 538, pushing the existing rows 538:620 down to 539:621.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("538:538").Insert()

$ws.Range("A538").Value = 9
$ws.Range("B538").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C538").Value = "Metropolitana"
$ws.Range("D538").Value = 45180
$ws.Range("E538").Value = 13
$ws.Range("F538").Value = 100112052
$ws.Range("G538").Value = "Albahaca"
$ws.Range("H538").Value = "Sin especificar"
$ws.Range("I538").Value = "Primera"
$ws.Range("J538").Value = 340
$ws.Range("K538").Value = 4500
$ws.Range("L538").Value = 5000
$ws.Range("M538").Value = 4750
$ws.Range("N538").Value = "$/paquete"
$ws.Range("O538").Value = "Región de Arica y Parinacota"
$ws.Range("P538").Value = 4750
$ws.Range("Q538").Value = 1
$ws.Range("R538").Value = "Hortaliza"
